$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the "Theme Party" activity to "Kolokium Zon Selatan" and
#    reset its merit figure back to 0 (final 2023 generation update).
$ws.Range("C40").Value = "Kolokium Zon Selatan"
$ws.Range("D40").Value = 0

# 2. Bring the "Penandaan Fail" category row (row 15) into line with the
#    other category rows (21, 27, 33): give C15 the same "blank merged
#    entry" formatting used by C33 (no wrap/shrink, plain default font,
#    same thin top/bottom border), then merge B15:C15 like the others.
$ws.Range("C33").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B15:C15").Merge()
